# Auto-generated edit script: update cryptos list with new prices/volumes
# and re-rank WrappedliquidstakedEther2.0 and RocketPoolETH blocks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. '1.000', '0.9997') are preserved as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = '29.321.34'
$ws.Range("E2").Value = '  +0.42%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.860.79'
$ws.Range("E3").Value = '  +0.22%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5: XRP
$ws.Range("D5").Value = '0.7047'
$ws.Range("E5").Value = '  +0.90%  '

# Row 6: BNB
$ws.Range("D6").Value = '238.35'
$ws.Range("E6").Value = '  +0.52%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.00%  '

# Row 8: Dogecoin
$ws.Range("D8").Value = '0.07904'
$ws.Range("E8").Value = '  +2.29%  '

# Row 9: Cardano
$ws.Range("D9").Value = '0.3049'
$ws.Range("E9").Value = '  +0.28%  '

# Row 10: Solana
$ws.Range("D10").Value = '24.88'
$ws.Range("E10").Value = '  +7.05%  '

# Row 11: WrappedEther
$ws.Range("D11").Value = '2.028.92'
$ws.Range("E11").Value = '  +11.11%  '

# Row 12: TRON
$ws.Range("D12").Value = '0.08187'
$ws.Range("E12").Value = '  +0.32%  '

# Row 13: Polkadot
$ws.Range("D13").Value = '5.230'
$ws.Range("E13").Value = '  +1.55%  '

# Row 14: Polygon
$ws.Range("D14").Value = '0.7182'
$ws.Range("E14").Value = '  +0.21%  '

# Row 15: Litecoin
$ws.Range("D15").Value = '89.69'
$ws.Range("E15").Value = '  +0.67%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '29.325.48'
$ws.Range("E16").Value = '  +0.41%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.378.75'
$ws.Range("E17").Value = '  +12.62%  '

# Row 18: Uniswap
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '5.837'
$ws.Range("E18").Value = '  +1.39%  '

# Row 19: ShibaInu
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007814'
$ws.Range("E19").Value = '  +1.28%  '

# Row 20: Avalanche
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  -0.53%  '

# Row 21: BitcoinCash
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '238.81'
$ws.Range("E21").Value = '  +0.93%  '

# Row 22: Dai
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  +0.07%  '

# Row 23: BinanceUSD
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '1.0000'
$ws.Range("E23").Value = '  -0.04%  '

# Row 24: Chainlink
$ws.Range("D24").Value = '7.577'
$ws.Range("E24").Value = '  +2.03%  '

# Row 25: Monero
$ws.Range("D25").Value = '162.98'
$ws.Range("E25").Value = '  +0.35%  '

# Row 26: Cosmos
$ws.Range("D26").Value = '8.917'
$ws.Range("E26").Value = '  -0.95%  '

# Row 27: Stellar
$ws.Range("D27").Value = '0.1431'
$ws.Range("E27").Value = '  -2.78%  '

# Row 28: EthereumClassic
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  +0.77%  '

# Row 29: LidoDAOToken
$ws.Range("D29").Value = '1.917'
$ws.Range("E29").Value = '  -6.53%  '

# Row 30: Toncoin
$ws.Range("D30").Value = '1.379'
$ws.Range("E30").Value = '  -2.76%  '

# Row 31: PancakeSwap
$ws.Range("D31").Value = '1.474'
$ws.Range("E31").Value = '  -0.47%  '

# Row 32: Filecoin
$ws.Range("D32").Value = '4.333'
$ws.Range("E32").Value = '  -2.15%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = '4.057'
$ws.Range("E33").Value = '  +0.83%  '

# Row 34: Hedera
$ws.Range("D34").Value = '0.05186'
$ws.Range("E34").Value = '  -0.34%  '

# Row 35: ARBITRUM
$ws.Range("D35").Value = '1.178'
$ws.Range("E35").Value = '  +1.11%  '

# Row 36: ImmutableX
$ws.Range("D36").Value = '0.7131'
$ws.Range("E36").Value = '  +0.90%  '

# Row 37: Frax
$ws.Range("E37").Value = '  +0.61%  '

# Row 38: HuobiToken
$ws.Range("D38").Value = '2.674'
$ws.Range("E38").Value = '  +0.13%  '

# Row 39: VeChain
$ws.Range("D39").Value = '0.01854'
$ws.Range("E39").Value = '  +0.56%  '

# Row 40: MXToken
$ws.Range("E40").Value = '  -1.22%  '

# Row 41: Maker
$ws.Range("D41").Value = '1.169.51'
$ws.Range("E41").Value = '  +3.05%  '

# Row 42: TrustWalletToken
$ws.Range("D42").Value = '0.9233'
$ws.Range("E42").Value = '  -1.17%  '

# Row 43: FraxShare
$ws.Range("D43").Value = '6.023'
$ws.Range("E43").Value = '  +2.23%  '

# Row 44: Aave
$ws.Range("D44").Value = '71.66'
$ws.Range("E44").Value = '  +1.56%  '

# Row 45: TheSandbox
$ws.Range("D45").Value = '0.4269'
$ws.Range("E45").Value = '  +0.03%  '

# Row 46: RocketPoolETH
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.224.96'
$ws.Range("E46").Value = '  +10.94%  '

# Row 47: PaxDollar
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.03%  '

# Row 48: Quant
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '101.79'
$ws.Range("E48").Value = '  -1.49%  '

# Row 49: Mantle
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.5331'
$ws.Range("E49").Value = '  -2.56%  '

# Row 50: RenderToken
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.761'
$ws.Range("E50").Value = '  -1.57%  '

# Row 51: EnergySwap
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.195'
$ws.Range("E51").Value = '  +0.53%  '

